$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> new F-column value }
$changes = @{
    "展览" = @{
        3  = 1880
        5  = 19
        6  = 848
        9  = 41
        16 = 4425
        17 = 15
        18 = 32
        19 = 481
        20 = 429
        21 = 9
        22 = 12
        23 = 1007
        24 = 1971
        25 = 371
        26 = 50
        27 = 28
        28 = 50
        29 = 2110
        30 = 76
        32 = 18
        33 = 149
        34 = 98
        35 = 34
        36 = 215
        37 = 30
    }
    "演出" = @{
        2 = 33
    }
    "全部类型" = @{
        3  = 1880
        5  = 19
        6  = 848
        9  = 41
        16 = 33
        17 = 4425
        18 = 15
        19 = 32
        20 = 481
        21 = 429
        22 = 9
        23 = 12
        24 = 1007
        25 = 1971
        26 = 371
        27 = 50
        28 = 28
        29 = 50
        30 = 2110
        31 = 76
        33 = 18
        34 = 149
        35 = 98
        36 = 34
        37 = 215
        38 = 30
    }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $changes[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowMap[$row]
    }
}
